$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enable iterative calculation with a small convergence delta (iterateDelta) ---
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# --- Update unit headers (row 2) ---
# Order matters for shared-string pool ordering: change C2 before B2
$ws.Range("C2").Value2 = "kT/GWh"
$ws.Range("B2").Value2 = "M€/GWhth"

# --- Rescale cost column (B) from k-units to M-units / MW to GW (divide by 1000) ---
$ws.Range("B3").Value2 = 35 / 1000
$ws.Range("B4").Value2 = 2.5 / 1000
$ws.Range("B5").Value2 = 9 / 1000
$ws.Range("B6").Value2 = 8.5 / 1000
$ws.Range("B8").Value2 = 73 / 1000

# --- Batteries (CO2 row) costs updated, and CO2 cost ($/t -> M$/kT) + red font highlight ---
$ws.Range("B9").Value2 = 45 / 1000
$ws.Range("C9").Value2 = 0
$ws.Range("B9").Font.Color = 255

# --- New row 10: unit label for the newly added CO2 cost column, same red styling ---
$ws.Range("B10").Value2 = "M€/kTCO2"
$ws.Range("B10").Font.Color = 255

# --- Update selection to match author's final cursor position ---
$null = $ws.Range("D11").Select()
